$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold plain text values (coin names, links, price
# strings such as "480.00"/"603.72", and padded percentage strings).
# Force the Text number format first so Excel does not silently
# reinterpret number-looking strings (losing trailing zeros / exact
# formatting, or turning "+1.39%" into a percentage value, etc.).

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "63.818.70"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "3.320.03"
$ws.Range("E3").Value = "  +5.26%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "603.72"
$ws.Range("E5").Value = "  +2.51%  "

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "142.49"
$ws.Range("E6").Value = "  +3.07%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "3.318.42"
$ws.Range("E8").Value = "  +5.21%  "

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +0.86%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.10%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.37%  "

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  +2.51%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.79%  "

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "34.92"
$ws.Range("E14").Value = "  +2.55%  "

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "3.872.05"
$ws.Range("E15").Value = "  +5.50%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.322.13"
$ws.Range("E17").Value = "  +5.37%  "

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "63.880.39"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("E19").Value = "  +3.51%  "

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "480.00"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "14.12"
$ws.Range("E21").Value = "  +1.53%  "

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +5.27%  "

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "8.18"
$ws.Range("E23").Value = "  +6.09%  "

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "13.71"
$ws.Range("E24").Value = "  +5.78%  "

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.72%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("B29:E29").NumberFormat = "@"
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  +3.53%  "

# Row 30
$ws.Range("B30:E30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "8.21"
$ws.Range("E30").Value = "  +3.57%  "

# Row 31
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").Value = "  +3.82%  "

# Row 32
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "28.84"
$ws.Range("E32").Value = "  +7.73%  "

# Row 34
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  +0.59%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.76%  "

# Row 36
$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  +5.04%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.21%  "

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  +6.18%  "

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0400"
$ws.Range("E39").Value = "  +3.78%  "

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "433.64"
$ws.Range("E40").Value = "  +4.13%  "

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "3.105.98"
$ws.Range("E41").Value = "  +5.19%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +9.81%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "8.34"
$ws.Range("E44").Value = "  +1.07%  "

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  +2.59%  "

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  +4.81%  "

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "37.27"
$ws.Range("E47").Value = "  +15.85%  "

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "26.36"
$ws.Range("E48").Value = "  +3.83%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("B50:E50").NumberFormat = "@"
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +3.04%  "

# Row 51
$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.114"
$ws.Range("E51").Value = "  +0.99%  "

